$wb = $excel.ActiveWorkbook

# --- Sheet 1 ---
$ws = $wb.Worksheets.Item(1)
$ws.Name = "summ45596795"
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = [double]"9775.768023736131"
$ws.Cells.Item(2, 3).Value = [double]"2.075684542651179e-255"
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = [double]"-1149.246552810493"
$ws.Cells.Item(3, 3).Value = [double]"1.053692263298837e-20"
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = [double]"-623.2866425669056"
$ws.Cells.Item(4, 3).Value = [double]"1.678400148541322e-16"
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = [double]"-134.3993086671417"
$ws.Cells.Item(5, 3).Value = [double]"0.05076352485662478"
$ws.Cells.Item(6, 1).Value = "Season[T.Spring]"
$ws.Cells.Item(6, 2).Value = [double]"10.45997317822565"
$ws.Cells.Item(6, 3).Value = [double]"0.8889591708601026"
$ws.Cells.Item(7, 1).Value = "Season[T.Summer]"
$ws.Cells.Item(7, 2).Value = [double]"-160.937130642773"
$ws.Cells.Item(7, 3).Value = [double]"0.1079085032676634"
$ws.Cells.Item(8, 1).Value = "Season[T.Winter]"
$ws.Cells.Item(8, 2).Value = [double]"-379.3116428226595"
$ws.Cells.Item(8, 3).Value = [double]"4.532566148135785e-08"
$ws.Cells.Item(9, 1).Value = "Country[T.France]"
$ws.Cells.Item(9, 2).Value = [double]"-1820.218281795273"
$ws.Cells.Item(9, 3).Value = [double]"1.229999153549695e-44"
$ws.Cells.Item(10, 1).Value = "Country[T.Germany]"
$ws.Cells.Item(10, 2).Value = [double]"-1103.354676352577"
$ws.Cells.Item(10, 3).Value = [double]"1.405422033159878e-17"
$ws.Cells.Item(11, 1).Value = "Country[T.Spain]"
$ws.Cells.Item(11, 2).Value = [double]"-1673.575758853573"
$ws.Cells.Item(11, 3).Value = [double]"7.752010702783338e-32"
$ws.Cells.Item(12, 1).Value = "HHSize"
$ws.Cells.Item(12, 2).Value = [double]"5.654694354318963"
$ws.Cells.Item(12, 3).Value = [double]"0.7430442749630252"
$ws.Cells.Item(13, 1).Value = "Sex"
$ws.Cells.Item(13, 2).Value = [double]"-1212.548588852925"
$ws.Cells.Item(13, 3).Value = [double]"2.126930980609474e-184"
$ws.Cells.Item(14, 1).Value = "Age"
$ws.Cells.Item(14, 2).Value = [double]"-18.33605363761876"
$ws.Cells.Item(14, 3).Value = [double]"1.623509060681682e-24"
$ws.Cells.Item(15, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(15, 2).Value = [double]"155.6255857099081"
$ws.Cells.Item(15, 3).Value = [double]"1.049910286596362e-47"
$ws.Cells.Item(16, 1).Value = "DistCenter_res"
$ws.Cells.Item(16, 2).Value = [double]"371.2212237734104"
$ws.Cells.Item(16, 3).Value = [double]"0"
$ws.Cells.Item(17, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(17, 2).Value = [double]"-0.02045634124758675"
$ws.Cells.Item(17, 3).Value = [double]"7.730069114163102e-09"
$ws.Cells.Item(18, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(18, 2).Value = [double]"-5.051033048071565e-06"
$ws.Cells.Item(18, 3).Value = [double]"0.04486641322755662"
$ws.Cells.Item(19, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(19, 2).Value = [double]"-12.58104534592443"
$ws.Cells.Item(19, 3).Value = [double]"2.225563830097688e-23"
$ws.Cells.Item(20, 1).Value = "street_length_res"
$ws.Cells.Item(20, 2).Value = [double]"9.62884810699396"
$ws.Cells.Item(20, 3).Value = [double]"5.648801283787504e-11"
$ws.Cells.Item(21, 1).Value = "LU_Comm_res"
$ws.Cells.Item(21, 2).Value = [double]"-3392.570948128269"
$ws.Cells.Item(21, 3).Value = [double]"2.507348979795321e-54"
$ws.Cells.Item(22, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(22, 2).Value = [double]"-1614.917563189271"
$ws.Cells.Item(22, 3).Value = [double]"9.270268287045025e-26"
$ws.Cells.Item(23, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(23, 2).Value = [double]"224.0085050308739"
$ws.Cells.Item(23, 3).Value = [double]"0.3265382861105424"

# --- Sheet 2 ---
$ws = $wb.Worksheets.Item(2)
$ws.Name = "summ54307975"
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = [double]"10084.45268680887"
$ws.Cells.Item(2, 3).Value = [double]"2.418049489708162e-272"
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = [double]"-1220.655073908245"
$ws.Cells.Item(3, 3).Value = [double]"2.969943534235572e-23"
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = [double]"-768.9834659548053"
$ws.Cells.Item(4, 3).Value = [double]"2.515240275701952e-24"
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = [double]"-216.4960858575028"
$ws.Cells.Item(5, 3).Value = [double]"0.001646990936878903"
$ws.Cells.Item(6, 1).Value = "Season[T.Spring]"
$ws.Cells.Item(6, 2).Value = [double]"46.15542337470166"
$ws.Cells.Item(6, 3).Value = [double]"0.5378333643046296"
$ws.Cells.Item(7, 1).Value = "Season[T.Summer]"
$ws.Cells.Item(7, 2).Value = [double]"-63.7995076917405"
$ws.Cells.Item(7, 3).Value = [double]"0.5235607303080329"
$ws.Cells.Item(8, 1).Value = "Season[T.Winter]"
$ws.Cells.Item(8, 2).Value = [double]"-322.8136031723171"
$ws.Cells.Item(8, 3).Value = [double]"3.397875223837344e-06"
$ws.Cells.Item(9, 1).Value = "Country[T.France]"
$ws.Cells.Item(9, 2).Value = [double]"-1969.655073275215"
$ws.Cells.Item(9, 3).Value = [double]"1.532385652302459e-51"
$ws.Cells.Item(10, 1).Value = "Country[T.Germany]"
$ws.Cells.Item(10, 2).Value = [double]"-1197.840937139137"
$ws.Cells.Item(10, 3).Value = [double]"2.646218811738014e-20"
$ws.Cells.Item(11, 1).Value = "Country[T.Spain]"
$ws.Cells.Item(11, 2).Value = [double]"-1859.42526858429"
$ws.Cells.Item(11, 3).Value = [double]"1.126854782929185e-38"
$ws.Cells.Item(12, 1).Value = "HHSize"
$ws.Cells.Item(12, 2).Value = [double]"22.70268462119851"
$ws.Cells.Item(12, 3).Value = [double]"0.1881338273136408"
$ws.Cells.Item(13, 1).Value = "Sex"
$ws.Cells.Item(13, 2).Value = [double]"-1214.722261944112"
$ws.Cells.Item(13, 3).Value = [double]"6.143975597863099e-185"
$ws.Cells.Item(14, 1).Value = "Age"
$ws.Cells.Item(14, 2).Value = [double]"-19.13193872545177"
$ws.Cells.Item(14, 3).Value = [double]"1.437681161146971e-26"
$ws.Cells.Item(15, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(15, 2).Value = [double]"145.8536822754413"
$ws.Cells.Item(15, 3).Value = [double]"5.890288919883858e-42"
$ws.Cells.Item(16, 1).Value = "DistCenter_res"
$ws.Cells.Item(16, 2).Value = [double]"371.3392358405202"
$ws.Cells.Item(16, 3).Value = [double]"0"
$ws.Cells.Item(17, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(17, 2).Value = [double]"-0.02039379688660448"
$ws.Cells.Item(17, 3).Value = [double]"8.474897979308623e-09"
$ws.Cells.Item(18, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(18, 2).Value = [double]"-5.173791942430197e-06"
$ws.Cells.Item(18, 3).Value = [double]"0.0405663844246937"
$ws.Cells.Item(19, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(19, 2).Value = [double]"-12.94698026957055"
$ws.Cells.Item(19, 3).Value = [double]"1.269846093848403e-24"
$ws.Cells.Item(20, 1).Value = "street_length_res"
$ws.Cells.Item(20, 2).Value = [double]"8.273635202976156"
$ws.Cells.Item(20, 3).Value = [double]"1.584680987793397e-08"
$ws.Cells.Item(21, 1).Value = "LU_Comm_res"
$ws.Cells.Item(21, 2).Value = [double]"-3277.59699635009"
$ws.Cells.Item(21, 3).Value = [double]"7.334161689313368e-51"
$ws.Cells.Item(22, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(22, 2).Value = [double]"-1535.443613753618"
$ws.Cells.Item(22, 3).Value = [double]"1.800372225674672e-23"
$ws.Cells.Item(23, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(23, 2).Value = [double]"260.9429002162485"
$ws.Cells.Item(23, 3).Value = [double]"0.2545298803714434"

# --- Sheet 3 ---
$ws = $wb.Worksheets.Item(3)
$ws.Name = "summ03859591"
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = [double]"9810.656300907242"
$ws.Cells.Item(2, 3).Value = [double]"9.502098459267502e-258"
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = [double]"-1299.510820445098"
$ws.Cells.Item(3, 3).Value = [double]"3.633634638154525e-26"
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = [double]"-750.0138678277885"
$ws.Cells.Item(4, 3).Value = [double]"3.010494790449694e-23"
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = [double]"-179.232543239235"
$ws.Cells.Item(5, 3).Value = [double]"0.009060175039805992"
$ws.Cells.Item(6, 1).Value = "Season[T.Spring]"
$ws.Cells.Item(6, 2).Value = [double]"20.91992141064203"
$ws.Cells.Item(6, 3).Value = [double]"0.779967269299312"
$ws.Cells.Item(7, 1).Value = "Season[T.Summer]"
$ws.Cells.Item(7, 2).Value = [double]"-165.715819366302"
$ws.Cells.Item(7, 3).Value = [double]"0.09585469913724035"
$ws.Cells.Item(8, 1).Value = "Season[T.Winter]"
$ws.Cells.Item(8, 2).Value = [double]"-356.4713903107764"
$ws.Cells.Item(8, 3).Value = [double]"2.74571015594336e-07"
$ws.Cells.Item(9, 1).Value = "Country[T.France]"
$ws.Cells.Item(9, 2).Value = [double]"-1769.043074018652"
$ws.Cells.Item(9, 3).Value = [double]"5.05680675265012e-42"
$ws.Cells.Item(10, 1).Value = "Country[T.Germany]"
$ws.Cells.Item(10, 2).Value = [double]"-1021.427981900328"
$ws.Cells.Item(10, 3).Value = [double]"3.135572559084717e-15"
$ws.Cells.Item(11, 1).Value = "Country[T.Spain]"
$ws.Cells.Item(11, 2).Value = [double]"-1654.186863428455"
$ws.Cells.Item(11, 3).Value = [double]"5.013704361112558e-31"
$ws.Cells.Item(12, 1).Value = "HHSize"
$ws.Cells.Item(12, 2).Value = [double]"33.69970257231529"
$ws.Cells.Item(12, 3).Value = [double]"0.05084582561078858"
$ws.Cells.Item(13, 1).Value = "Sex"
$ws.Cells.Item(13, 2).Value = [double]"-1168.055571999082"
$ws.Cells.Item(13, 3).Value = [double]"1.144343773519352e-171"
$ws.Cells.Item(14, 1).Value = "Age"
$ws.Cells.Item(14, 2).Value = [double]"-19.6034660712961"
$ws.Cells.Item(14, 3).Value = [double]"8.705619503240654e-28"
$ws.Cells.Item(15, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(15, 2).Value = [double]"153.5956503671911"
$ws.Cells.Item(15, 3).Value = [double]"1.618324342346258e-46"
$ws.Cells.Item(16, 1).Value = "DistCenter_res"
$ws.Cells.Item(16, 2).Value = [double]"370.5936190625175"
$ws.Cells.Item(16, 3).Value = [double]"0"
$ws.Cells.Item(17, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(17, 2).Value = [double]"-0.01990304236401001"
$ws.Cells.Item(17, 3).Value = [double]"1.898025428270954e-08"
$ws.Cells.Item(18, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(18, 2).Value = [double]"-4.437994255433712e-06"
$ws.Cells.Item(18, 3).Value = [double]"0.07670596331590561"
$ws.Cells.Item(19, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(19, 2).Value = [double]"-12.90691314848075"
$ws.Cells.Item(19, 3).Value = [double]"1.521440319074497e-24"
$ws.Cells.Item(20, 1).Value = "street_length_res"
$ws.Cells.Item(20, 2).Value = [double]"8.822462824644628"
$ws.Cells.Item(20, 3).Value = [double]"1.638688803367205e-09"
$ws.Cells.Item(21, 1).Value = "LU_Comm_res"
$ws.Cells.Item(21, 2).Value = [double]"-3363.184067106884"
$ws.Cells.Item(21, 3).Value = [double]"1.902123423106526e-53"
$ws.Cells.Item(22, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(22, 2).Value = [double]"-1564.289265293549"
$ws.Cells.Item(22, 3).Value = [double]"1.968269988421878e-24"
$ws.Cells.Item(23, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(23, 2).Value = [double]"250.2938625394086"
$ws.Cells.Item(23, 3).Value = [double]"0.2717900232120699"

# --- Sheet 4 ---
$ws = $wb.Worksheets.Item(4)
$ws.Name = "summ11665252"
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = [double]"9867.831762570202"
$ws.Cells.Item(2, 3).Value = [double]"2.700160097709771e-262"
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = [double]"-1255.178335248698"
$ws.Cells.Item(3, 3).Value = [double]"2.602411378413879e-24"
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = [double]"-694.0942574058545"
$ws.Cells.Item(4, 3).Value = [double]"3.651574927041322e-20"
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = [double]"-160.5066322722692"
$ws.Cells.Item(5, 3).Value = [double]"0.01934537350374139"
$ws.Cells.Item(6, 1).Value = "Season[T.Spring]"
$ws.Cells.Item(6, 2).Value = [double]"-7.526012882906976"
$ws.Cells.Item(6, 3).Value = [double]"0.9198157371583289"
$ws.Cells.Item(7, 1).Value = "Season[T.Summer]"
$ws.Cells.Item(7, 2).Value = [double]"-159.3882480229053"
$ws.Cells.Item(7, 3).Value = [double]"0.1100486465786124"
$ws.Cells.Item(8, 1).Value = "Season[T.Winter]"
$ws.Cells.Item(8, 2).Value = [double]"-340.4173760107034"
$ws.Cells.Item(8, 3).Value = [double]"8.727548038102442e-07"
$ws.Cells.Item(9, 1).Value = "Country[T.France]"
$ws.Cells.Item(9, 2).Value = [double]"-1865.984705936577"
$ws.Cells.Item(9, 3).Value = [double]"7.490065597298348e-47"
$ws.Cells.Item(10, 1).Value = "Country[T.Germany]"
$ws.Cells.Item(10, 2).Value = [double]"-1058.001924489246"
$ws.Cells.Item(10, 3).Value = [double]"2.582744357028173e-16"
$ws.Cells.Item(11, 1).Value = "Country[T.Spain]"
$ws.Cells.Item(11, 2).Value = [double]"-1610.733959985766"
$ws.Cells.Item(11, 3).Value = [double]"1.207175287523071e-29"
$ws.Cells.Item(12, 1).Value = "HHSize"
$ws.Cells.Item(12, 2).Value = [double]"-3.464251889374662"
$ws.Cells.Item(12, 3).Value = [double]"0.8406949114826485"
$ws.Cells.Item(13, 1).Value = "Sex"
$ws.Cells.Item(13, 2).Value = [double]"-1173.128026072646"
$ws.Cells.Item(13, 3).Value = [double]"1.588587114347209e-173"
$ws.Cells.Item(14, 1).Value = "Age"
$ws.Cells.Item(14, 2).Value = [double]"-19.36699068211325"
$ws.Cells.Item(14, 3).Value = [double]"3.028719302456937e-27"
$ws.Cells.Item(15, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(15, 2).Value = [double]"144.4709282631736"
$ws.Cells.Item(15, 3).Value = [double]"2.109143788196584e-41"
$ws.Cells.Item(16, 1).Value = "DistCenter_res"
$ws.Cells.Item(16, 2).Value = [double]"372.3934161142266"
$ws.Cells.Item(16, 3).Value = [double]"0"
$ws.Cells.Item(17, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(17, 2).Value = [double]"-0.02096280358801943"
$ws.Cells.Item(17, 3).Value = [double]"3.206139197778252e-09"
$ws.Cells.Item(18, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(18, 2).Value = [double]"-5.587592028862135e-06"
$ws.Cells.Item(18, 3).Value = [double]"0.02566781545097941"
$ws.Cells.Item(19, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(19, 2).Value = [double]"-13.43656897642266"
$ws.Cells.Item(19, 3).Value = [double]"1.40579262099312e-26"
$ws.Cells.Item(20, 1).Value = "street_length_res"
$ws.Cells.Item(20, 2).Value = [double]"8.09969822378133"
$ws.Cells.Item(20, 3).Value = [double]"2.853857576504282e-08"
$ws.Cells.Item(21, 1).Value = "LU_Comm_res"
$ws.Cells.Item(21, 2).Value = [double]"-3036.36645178668"
$ws.Cells.Item(21, 3).Value = [double]"3.096086030510595e-44"
$ws.Cells.Item(22, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(22, 2).Value = [double]"-1256.992390702074"
$ws.Cells.Item(22, 3).Value = [double]"2.257042097185105e-16"
$ws.Cells.Item(23, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(23, 2).Value = [double]"163.5548847721138"
$ws.Cells.Item(23, 3).Value = [double]"0.4751108679742072"

# --- Sheet 5 ---
$ws = $wb.Worksheets.Item(5)
$ws.Name = "summ19378542"
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = [double]"9812.068427009559"
$ws.Cells.Item(2, 3).Value = [double]"3.011539093636567e-257"
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = [double]"-1080.598206180976"
$ws.Cells.Item(3, 3).Value = [double]"1.575862955396212e-18"
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = [double]"-692.7113012861928"
$ws.Cells.Item(4, 3).Value = [double]"6.563545027263914e-20"
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = [double]"-157.1149017549513"
$ws.Cells.Item(5, 3).Value = [double]"0.02272180091163794"
$ws.Cells.Item(6, 1).Value = "Season[T.Spring]"
$ws.Cells.Item(6, 2).Value = [double]"-12.31251023839087"
$ws.Cells.Item(6, 3).Value = [double]"0.8697109164389478"
$ws.Cells.Item(7, 1).Value = "Season[T.Summer]"
$ws.Cells.Item(7, 2).Value = [double]"-105.9944254592187"
$ws.Cells.Item(7, 3).Value = [double]"0.2890968227686593"
$ws.Cells.Item(8, 1).Value = "Season[T.Winter]"
$ws.Cells.Item(8, 2).Value = [double]"-353.6343057007252"
$ws.Cells.Item(8, 3).Value = [double]"3.576751438989741e-07"
$ws.Cells.Item(9, 1).Value = "Country[T.France]"
$ws.Cells.Item(9, 2).Value = [double]"-1740.494731663622"
$ws.Cells.Item(9, 3).Value = [double]"2.254600274107986e-40"
$ws.Cells.Item(10, 1).Value = "Country[T.Germany]"
$ws.Cells.Item(10, 2).Value = [double]"-967.0984269623314"
$ws.Cells.Item(10, 3).Value = [double]"1.080931672980224e-13"
$ws.Cells.Item(11, 1).Value = "Country[T.Spain]"
$ws.Cells.Item(11, 2).Value = [double]"-1527.958804779192"
$ws.Cells.Item(11, 3).Value = [double]"1.746628199122019e-26"
$ws.Cells.Item(12, 1).Value = "HHSize"
$ws.Cells.Item(12, 2).Value = [double]"18.17255864318427"
$ws.Cells.Item(12, 3).Value = [double]"0.2927196066133477"
$ws.Cells.Item(13, 1).Value = "Sex"
$ws.Cells.Item(13, 2).Value = [double]"-1193.349920815175"
$ws.Cells.Item(13, 3).Value = [double]"5.176727610879944e-178"
$ws.Cells.Item(14, 1).Value = "Age"
$ws.Cells.Item(14, 2).Value = [double]"-19.41952973456577"
$ws.Cells.Item(14, 3).Value = [double]"3.259765164388117e-27"
$ws.Cells.Item(15, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(15, 2).Value = [double]"152.0432126107808"
$ws.Cells.Item(15, 3).Value = [double]"6.527373501356663e-45"
$ws.Cells.Item(16, 1).Value = "DistCenter_res"
$ws.Cells.Item(16, 2).Value = [double]"371.5129097908444"
$ws.Cells.Item(16, 3).Value = [double]"0"
$ws.Cells.Item(17, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(17, 2).Value = [double]"-0.01923474643553832"
$ws.Cells.Item(17, 3).Value = [double]"6.448769628043414e-08"
$ws.Cells.Item(18, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(18, 2).Value = [double]"-5.435439166229918e-06"
$ws.Cells.Item(18, 3).Value = [double]"0.03124591170610608"
$ws.Cells.Item(19, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(19, 2).Value = [double]"-13.37779950492716"
$ws.Cells.Item(19, 3).Value = [double]"3.19492004243805e-26"
$ws.Cells.Item(20, 1).Value = "street_length_res"
$ws.Cells.Item(20, 2).Value = [double]"8.635242750590933"
$ws.Cells.Item(20, 3).Value = [double]"4.278564494760932e-09"
$ws.Cells.Item(21, 1).Value = "LU_Comm_res"
$ws.Cells.Item(21, 2).Value = [double]"-3447.546712336319"
$ws.Cells.Item(21, 3).Value = [double]"5.642030352184246e-56"
$ws.Cells.Item(22, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(22, 2).Value = [double]"-1552.841084893899"
$ws.Cells.Item(22, 3).Value = [double]"7.735057180364184e-24"
$ws.Cells.Item(23, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(23, 2).Value = [double]"367.0187008583057"
$ws.Cells.Item(23, 3).Value = [double]"0.1099951521831636"

# --- Sheet 6 ---
$ws = $wb.Worksheets.Item(6)
$ws.Name = "summ30904350"
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = [double]"9749.517175811401"
$ws.Cells.Item(2, 3).Value = [double]"2.101185651845246e-255"
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = [double]"-1151.633404813422"
$ws.Cells.Item(3, 3).Value = [double]"6.577522111055699e-21"
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = [double]"-671.6111138056101"
$ws.Cells.Item(4, 3).Value = [double]"5.126319516849032e-19"
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = [double]"-149.9925749309848"
$ws.Cells.Item(5, 3).Value = [double]"0.02877266768261105"
$ws.Cells.Item(6, 1).Value = "Season[T.Spring]"
$ws.Cells.Item(6, 2).Value = [double]"-8.043810160003602"
$ws.Cells.Item(6, 3).Value = [double]"0.9141764116821613"
$ws.Cells.Item(7, 1).Value = "Season[T.Summer]"
$ws.Cells.Item(7, 2).Value = [double]"-222.4168649942185"
$ws.Cells.Item(7, 3).Value = [double]"0.02559313527892125"
$ws.Cells.Item(8, 1).Value = "Season[T.Winter]"
$ws.Cells.Item(8, 2).Value = [double]"-353.102194974777"
$ws.Cells.Item(8, 3).Value = [double]"3.341079785042222e-07"
$ws.Cells.Item(9, 1).Value = "Country[T.France]"
$ws.Cells.Item(9, 2).Value = [double]"-1879.477967332911"
$ws.Cells.Item(9, 3).Value = [double]"2.025647889044218e-47"
$ws.Cells.Item(10, 1).Value = "Country[T.Germany]"
$ws.Cells.Item(10, 2).Value = [double]"-1107.455506302689"
$ws.Cells.Item(10, 3).Value = [double]"1.029198844521303e-17"
$ws.Cells.Item(11, 1).Value = "Country[T.Spain]"
$ws.Cells.Item(11, 2).Value = [double]"-1675.56063743941"
$ws.Cells.Item(11, 3).Value = [double]"6.475133761870066e-32"
$ws.Cells.Item(12, 1).Value = "HHSize"
$ws.Cells.Item(12, 2).Value = [double]"21.32047472476322"
$ws.Cells.Item(12, 3).Value = [double]"0.2146022331696382"
$ws.Cells.Item(13, 1).Value = "Sex"
$ws.Cells.Item(13, 2).Value = [double]"-1133.134218223814"
$ws.Cells.Item(13, 3).Value = [double]"2.703927598455871e-162"
$ws.Cells.Item(14, 1).Value = "Age"
$ws.Cells.Item(14, 2).Value = [double]"-19.31233964987491"
$ws.Cells.Item(14, 3).Value = [double]"3.258467924312484e-27"
$ws.Cells.Item(15, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(15, 2).Value = [double]"151.0896185263566"
$ws.Cells.Item(15, 3).Value = [double]"4.488856828734076e-45"
$ws.Cells.Item(16, 1).Value = "DistCenter_res"
$ws.Cells.Item(16, 2).Value = [double]"373.9804834686414"
$ws.Cells.Item(16, 3).Value = [double]"0"
$ws.Cells.Item(17, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(17, 2).Value = [double]"-0.01824950691367757"
$ws.Cells.Item(17, 3).Value = [double]"2.607590560566457e-07"
$ws.Cells.Item(18, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(18, 2).Value = [double]"-5.214804445807742e-06"
$ws.Cells.Item(18, 3).Value = [double]"0.0374555850787077"
$ws.Cells.Item(19, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(19, 2).Value = [double]"-12.94216514356365"
$ws.Cells.Item(19, 3).Value = [double]"7.909145571035854e-25"
$ws.Cells.Item(20, 1).Value = "street_length_res"
$ws.Cells.Item(20, 2).Value = [double]"9.152292299854945"
$ws.Cells.Item(20, 3).Value = [double]"3.866060388865476e-10"
$ws.Cells.Item(21, 1).Value = "LU_Comm_res"
$ws.Cells.Item(21, 2).Value = [double]"-3229.150141566677"
$ws.Cells.Item(21, 3).Value = [double]"1.303573395608644e-49"
$ws.Cells.Item(22, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(22, 2).Value = [double]"-1564.852409822101"
$ws.Cells.Item(22, 3).Value = [double]"1.638634992818792e-24"
$ws.Cells.Item(23, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(23, 2).Value = [double]"308.820646832495"
$ws.Cells.Item(23, 3).Value = [double]"0.1747878408857753"

# --- Sheet 7 ---
$ws = $wb.Worksheets.Item(7)
$ws.Name = "summ42791849"
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = [double]"9638.793635027199"
$ws.Cells.Item(2, 3).Value = [double]"1.666789996266829e-246"
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = [double]"-1246.798284212027"
$ws.Cells.Item(3, 3).Value = [double]"4.070774587160402e-24"
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = [double]"-645.3365041640509"
$ws.Cells.Item(4, 3).Value = [double]"1.808421498098663e-17"
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = [double]"-144.7512231265257"
$ws.Cells.Item(5, 3).Value = [double]"0.03608969798984454"
$ws.Cells.Item(6, 1).Value = "Season[T.Spring]"
$ws.Cells.Item(6, 2).Value = [double]"28.57574589863758"
$ws.Cells.Item(6, 3).Value = [double]"0.7036782396262573"
$ws.Cells.Item(7, 1).Value = "Season[T.Summer]"
$ws.Cells.Item(7, 2).Value = [double]"-135.7074708267997"
$ws.Cells.Item(7, 3).Value = [double]"0.1752365655428683"
$ws.Cells.Item(8, 1).Value = "Season[T.Winter]"
$ws.Cells.Item(8, 2).Value = [double]"-320.4039482797058"
$ws.Cells.Item(8, 3).Value = [double]"4.081420801204785e-06"
$ws.Cells.Item(9, 1).Value = "Country[T.France]"
$ws.Cells.Item(9, 2).Value = [double]"-1769.503425507719"
$ws.Cells.Item(9, 3).Value = [double]"9.613185301496668e-42"
$ws.Cells.Item(10, 1).Value = "Country[T.Germany]"
$ws.Cells.Item(10, 2).Value = [double]"-1048.1142337981"
$ws.Cells.Item(10, 3).Value = [double]"7.229817709080291e-16"
$ws.Cells.Item(11, 1).Value = "Country[T.Spain]"
$ws.Cells.Item(11, 2).Value = [double]"-1653.256251746171"
$ws.Cells.Item(11, 3).Value = [double]"9.124411538606253e-31"
$ws.Cells.Item(12, 1).Value = "HHSize"
$ws.Cells.Item(12, 2).Value = [double]"22.65948334125458"
$ws.Cells.Item(12, 3).Value = [double]"0.1890737626824182"
$ws.Cells.Item(13, 1).Value = "Sex"
$ws.Cells.Item(13, 2).Value = [double]"-1189.965234962869"
$ws.Cells.Item(13, 3).Value = [double]"5.543792270926348e-177"
$ws.Cells.Item(14, 1).Value = "Age"
$ws.Cells.Item(14, 2).Value = [double]"-18.95881922539256"
$ws.Cells.Item(14, 3).Value = [double]"4.524990469961667e-26"
$ws.Cells.Item(15, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(15, 2).Value = [double]"150.7940073345918"
$ws.Cells.Item(15, 3).Value = [double]"9.763510289872004e-45"
$ws.Cells.Item(16, 1).Value = "DistCenter_res"
$ws.Cells.Item(16, 2).Value = [double]"374.3499445632093"
$ws.Cells.Item(16, 3).Value = [double]"0"
$ws.Cells.Item(17, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(17, 2).Value = [double]"-0.01854490504382837"
$ws.Cells.Item(17, 3).Value = [double]"1.713341054769107e-07"
$ws.Cells.Item(18, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(18, 2).Value = [double]"-4.389224883827492e-06"
$ws.Cells.Item(18, 3).Value = [double]"0.08030421676851732"
$ws.Cells.Item(19, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(19, 2).Value = [double]"-12.40313029529523"
$ws.Cells.Item(19, 3).Value = [double]"9.689081427247769e-23"
$ws.Cells.Item(20, 1).Value = "street_length_res"
$ws.Cells.Item(20, 2).Value = [double]"10.15229501496844"
$ws.Cells.Item(20, 3).Value = [double]"7.199738049593092e-12"
$ws.Cells.Item(21, 1).Value = "LU_Comm_res"
$ws.Cells.Item(21, 2).Value = [double]"-3526.010321790459"
$ws.Cells.Item(21, 3).Value = [double]"3.443701495524306e-58"
$ws.Cells.Item(22, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(22, 2).Value = [double]"-1690.400129960631"
$ws.Cells.Item(22, 3).Value = [double]"4.251793244370173e-28"
$ws.Cells.Item(23, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(23, 2).Value = [double]"296.0509542158907"
$ws.Cells.Item(23, 3).Value = [double]"0.1957251127081746"

# --- Sheet 8 ---
$ws = $wb.Worksheets.Item(8)
$ws.Name = "summ53533809"
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = [double]"9700.307453877882"
$ws.Cells.Item(2, 3).Value = [double]"3.078028267550654e-253"
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = [double]"-1234.311774351828"
$ws.Cells.Item(3, 3).Value = [double]"6.367252710575762e-24"
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = [double]"-668.1303039786387"
$ws.Cells.Item(4, 3).Value = [double]"7.234362267155637e-19"
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = [double]"-162.972237097095"
$ws.Cells.Item(5, 3).Value = [double]"0.01739521594649765"
$ws.Cells.Item(6, 1).Value = "Season[T.Spring]"
$ws.Cells.Item(6, 2).Value = [double]"13.70071812215403"
$ws.Cells.Item(6, 3).Value = [double]"0.854538922692719"
$ws.Cells.Item(7, 1).Value = "Season[T.Summer]"
$ws.Cells.Item(7, 2).Value = [double]"-135.9840534528213"
$ws.Cells.Item(7, 3).Value = [double]"0.1727930116416576"
$ws.Cells.Item(8, 1).Value = "Season[T.Winter]"
$ws.Cells.Item(8, 2).Value = [double]"-381.9833153067084"
$ws.Cells.Item(8, 3).Value = [double]"3.474868776301038e-08"
$ws.Cells.Item(9, 1).Value = "Country[T.France]"
$ws.Cells.Item(9, 2).Value = [double]"-1831.097734193478"
$ws.Cells.Item(9, 3).Value = [double]"3.117273645707762e-45"
$ws.Cells.Item(10, 1).Value = "Country[T.Germany]"
$ws.Cells.Item(10, 2).Value = [double]"-1155.306874709391"
$ws.Cells.Item(10, 3).Value = [double]"3.388457576697973e-19"
$ws.Cells.Item(11, 1).Value = "Country[T.Spain]"
$ws.Cells.Item(11, 2).Value = [double]"-1651.22591644504"
$ws.Cells.Item(11, 3).Value = [double]"3.720315569660791e-31"
$ws.Cells.Item(12, 1).Value = "HHSize"
$ws.Cells.Item(12, 2).Value = [double]"23.34450706647808"
$ws.Cells.Item(12, 3).Value = [double]"0.1753284305191965"
$ws.Cells.Item(13, 1).Value = "Sex"
$ws.Cells.Item(13, 2).Value = [double]"-1193.743034868023"
$ws.Cells.Item(13, 3).Value = [double]"1.65034989581659e-179"
$ws.Cells.Item(14, 1).Value = "Age"
$ws.Cells.Item(14, 2).Value = [double]"-18.94776294860595"
$ws.Cells.Item(14, 3).Value = [double]"3.952838914633366e-26"
$ws.Cells.Item(15, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(15, 2).Value = [double]"154.4210840431857"
$ws.Cells.Item(15, 3).Value = [double]"7.805802350085417e-47"
$ws.Cells.Item(16, 1).Value = "DistCenter_res"
$ws.Cells.Item(16, 2).Value = [double]"368.2693116028556"
$ws.Cells.Item(16, 3).Value = [double]"0"
$ws.Cells.Item(17, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(17, 2).Value = [double]"-0.01798575089460603"
$ws.Cells.Item(17, 3).Value = [double]"3.665943190854669e-07"
$ws.Cells.Item(18, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(18, 2).Value = [double]"-5.250940251046764e-06"
$ws.Cells.Item(18, 3).Value = [double]"0.03726173944122765"
$ws.Cells.Item(19, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(19, 2).Value = [double]"-13.82939925421223"
$ws.Cells.Item(19, 3).Value = [double]"5.346085944331097e-28"
$ws.Cells.Item(20, 1).Value = "street_length_res"
$ws.Cells.Item(20, 2).Value = [double]"10.34552645804578"
$ws.Cells.Item(20, 3).Value = [double]"1.553957155789994e-12"
$ws.Cells.Item(21, 1).Value = "LU_Comm_res"
$ws.Cells.Item(21, 2).Value = [double]"-3122.157061640347"
$ws.Cells.Item(21, 3).Value = [double]"2.215386927547556e-46"
$ws.Cells.Item(22, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(22, 2).Value = [double]"-1544.556863745444"
$ws.Cells.Item(22, 3).Value = [double]"8.382363635340312e-24"
$ws.Cells.Item(23, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(23, 2).Value = [double]"314.6325349596277"
$ws.Cells.Item(23, 3).Value = [double]"0.1669210739593829"

# --- Sheet 9 ---
$ws = $wb.Worksheets.Item(9)
$ws.Name = "summ03543543"
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = [double]"9897.521399287492"
$ws.Cells.Item(2, 3).Value = [double]"1.730340702658624e-261"
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = [double]"-1268.730010301057"
$ws.Cells.Item(3, 3).Value = [double]"4.350102217717693e-25"
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = [double]"-737.4432501737051"
$ws.Cells.Item(4, 3).Value = [double]"1.915847338526868e-22"
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = [double]"-204.1382801955077"
$ws.Cells.Item(5, 3).Value = [double]"0.002989573503644179"
$ws.Cells.Item(6, 1).Value = "Season[T.Spring]"
$ws.Cells.Item(6, 2).Value = [double]"11.11338332424071"
$ws.Cells.Item(6, 3).Value = [double]"0.8819594814401001"
$ws.Cells.Item(7, 1).Value = "Season[T.Summer]"
$ws.Cells.Item(7, 2).Value = [double]"-133.561507571111"
$ws.Cells.Item(7, 3).Value = [double]"0.1815911292955126"
$ws.Cells.Item(8, 1).Value = "Season[T.Winter]"
$ws.Cells.Item(8, 2).Value = [double]"-319.4633797298584"
$ws.Cells.Item(8, 3).Value = [double]"4.119032629437846e-06"
$ws.Cells.Item(9, 1).Value = "Country[T.France]"
$ws.Cells.Item(9, 2).Value = [double]"-1837.531779716148"
$ws.Cells.Item(9, 3).Value = [double]"5.34847982837994e-45"
$ws.Cells.Item(10, 1).Value = "Country[T.Germany]"
$ws.Cells.Item(10, 2).Value = [double]"-1094.701625219921"
$ws.Cells.Item(10, 3).Value = [double]"3.528346033365683e-17"
$ws.Cells.Item(11, 1).Value = "Country[T.Spain]"
$ws.Cells.Item(11, 2).Value = [double]"-1661.237014913399"
$ws.Cells.Item(11, 3).Value = [double]"3.821261262479006e-31"
$ws.Cells.Item(12, 1).Value = "HHSize"
$ws.Cells.Item(12, 2).Value = [double]"26.14564076248128"
$ws.Cells.Item(12, 3).Value = [double]"0.1293282139395703"
$ws.Cells.Item(13, 1).Value = "Sex"
$ws.Cells.Item(13, 2).Value = [double]"-1192.317778636344"
$ws.Cells.Item(13, 3).Value = [double]"1.426781377681877e-178"
$ws.Cells.Item(14, 1).Value = "Age"
$ws.Cells.Item(14, 2).Value = [double]"-19.59677379015218"
$ws.Cells.Item(14, 3).Value = [double]"9.177020463634362e-28"
$ws.Cells.Item(15, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(15, 2).Value = [double]"142.5541296986971"
$ws.Cells.Item(15, 3).Value = [double]"2.527901363401935e-40"
$ws.Cells.Item(16, 1).Value = "DistCenter_res"
$ws.Cells.Item(16, 2).Value = [double]"368.2423309287404"
$ws.Cells.Item(16, 3).Value = [double]"0"
$ws.Cells.Item(17, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(17, 2).Value = [double]"-0.0214248415372899"
$ws.Cells.Item(17, 3).Value = [double]"1.419212204588613e-09"
$ws.Cells.Item(18, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(18, 2).Value = [double]"-4.219812604212756e-06"
$ws.Cells.Item(18, 3).Value = [double]"0.09572439022910642"
$ws.Cells.Item(19, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(19, 2).Value = [double]"-12.89109992606573"
$ws.Cells.Item(19, 3).Value = [double]"1.597923495002545e-24"
$ws.Cells.Item(20, 1).Value = "street_length_res"
$ws.Cells.Item(20, 2).Value = [double]"8.766004540192931"
$ws.Cells.Item(20, 3).Value = [double]"2.255317349561811e-09"
$ws.Cells.Item(21, 1).Value = "LU_Comm_res"
$ws.Cells.Item(21, 2).Value = [double]"-3186.214823012045"
$ws.Cells.Item(21, 3).Value = [double]"3.31674836473813e-48"
$ws.Cells.Item(22, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(22, 2).Value = [double]"-1541.444988562078"
$ws.Cells.Item(22, 3).Value = [double]"1.161500613952264e-23"
$ws.Cells.Item(23, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(23, 2).Value = [double]"338.0450281111362"
$ws.Cells.Item(23, 3).Value = [double]"0.1379067979944381"

Write-Host "done"